# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces the two kinds of changes captured in the target diff:
#
#   1. Three tables (on the slides that hold a table) get their table
#      style switched from the custom style
#      {28A34EEA-2248-4B0C-876F-9FC12AC68CE6} to the built-in style
#      {F1A2B336-7303-4033-BCB7-023A6B306BD8}.
#
#   2. The presentation's theme color palette is switched from the
#      "Integral / Red Violet" palette to the standard "Office" palette
#      (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), which is what the
#      XML diff shows as the new content of ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table in the deck -----------------------------
$newTableStyleId = "{F1A2B336-7303-4033-BCB7-023A6B306BD8}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the theme color scheme to the "Office" palette ----------
# ThemeColorScheme index order: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1
# 6=accent2 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
$officeColors = @(
    0,         # dk1      #000000
    16777215,  # lt1      #FFFFFF
    6968388,   # dk2      #44546A
    15132391,  # lt2      #E7E6E6
    13998939,  # accent1  #5B9BD5
    3243501,   # accent2  #ED7D31
    10855845,  # accent3  #A5A5A5
    49407,     # accent4  #FFC000
    12874308,  # accent5  #4472C4
    4697456,   # accent6  #70AD47
    12673797,  # hlink    #0563C1
    7491477    # folHlink #954F72
)

$slide = $p.Slides.Item(1)
$colorScheme = $slide.ThemeColorScheme
for ($k = 1; $k -le $colorScheme.Count; $k++) {
    $colorScheme.Item($k).RGB = $officeColors[$k - 1]
}
